# Move the diagram shapes on slide 2 (the "model_diagram" task slide).
# All shapes are translated by the same vector (dx=419819 EMU, dy=1104182 EMU)
# except the three connectors that are glued to "Cylinder 21" (id 22) and
# re-route differently (one of them also grows slightly taller).
#
# NOTE: this COM-interop engine does not support PowerShell named
# parameters / default parameter values reliably, so Set-ShapeRect takes
# all five arguments positionally; pass -1 for CxEmu/CyEmu to leave the
# shape's current width/height untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$EMU = 12700.0

function Set-ShapeRect($Name, $XEmu, $YEmu, $CxEmu, $CyEmu) {
    $shp = $s.Shapes.Item($Name)
    $shp.Left = $XEmu / $EMU
    $shp.Top = $YEmu / $EMU
    if ($CxEmu -ge 0) {
        $shp.Width = $CxEmu / $EMU
    }
    if ($CyEmu -ge 0) {
        $shp.Height = $CyEmu / $EMU
    }
}

# Rectangle 58 (id 59)
Set-ShapeRect "Rectangle 58" 4271513 1220967 -1 -1
# Rectangle 48 (id 49)
Set-ShapeRect "Rectangle 48" 419819 1220967 -1 -1
# Rectangle 3 (id 4)
Set-ShapeRect "Rectangle 3" 544816 1697140 -1 -1
# Rectangle 4 (id 5)
Set-ShapeRect "Rectangle 4" 544815 2792694 -1 -1
# Rectangle 5 (id 6)
Set-ShapeRect "Rectangle 5" 544815 3888248 -1 -1
# Straight Arrow Connector 7 (id 8)
Set-ShapeRect "Straight Arrow Connector 7" 2244307 2172813 -1 -1
# Straight Arrow Connector 10 (id 11)
Set-ShapeRect "Straight Arrow Connector 10" 2244306 3253752 -1 -1
# Straight Arrow Connector 13 (id 14) - also grows slightly taller
Set-ShapeRect "Straight Arrow Connector 13" 2244306 4336367 767839 27554
# Cylinder 20 (id 21)
Set-ShapeRect "Cylinder 20" 3012145 1631505 -1 -1
# Cylinder 19 (id 20)
Set-ShapeRect "Cylinder 19" 3012145 2712444 -1 -1
# Cylinder 21 (id 22)
Set-ShapeRect "Cylinder 21" 3012145 3795059 -1 -1
# TextBox 49 (id 50)
Set-ShapeRect "TextBox 49" 665630 4977468 -1 -1
# Straight Arrow Connector 50 (id 51)
Set-ShapeRect "Straight Arrow Connector 50" 3848907 2172813 -1 -1
# Straight Arrow Connector 53 (id 54)
Set-ShapeRect "Straight Arrow Connector 53" 3848907 3253752 -1 -1
# Straight Arrow Connector 55 (id 56)
Set-ShapeRect "Straight Arrow Connector 55" 3848907 3268367 -1 -1
# TextBox 59 (id 60)
Set-ShapeRect "TextBox 59" 4271513 5254467 -1 -1
# Rectangle 60 (id 61)
Set-ShapeRect "Rectangle 60" 8195006 1220967 -1 -1
# TextBox 61 (id 62)
Set-ShapeRect "TextBox 61" 8401996 5254467 -1 -1
